$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 328.9565
$ws.Range("I39").Value = 274.5
$ws.Range("J39").Value = 334.14285
$ws.Range("K39").Value = 823.5
$ws.Range("L39").Value = 1002.42855
$ws.Range("M39").Value = -527.5
$ws.Range("N39").Value = -1594.42855

$ws.Range("H52").Value = 734.8823
$ws.Range("J52").Value = 9999
$ws.Range("L52").Value = 29997
$ws.Range("N52").Value = -30317

$ws.Range("H137").Value = 5066.619
$ws.Range("I137").Value = 2760.7917
$ws.Range("K137").Value = 8282.375100000001
$ws.Range("M137").Value = -5732.375100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 13188361
$ws.Range("I61").Value = 14709155
$ws.Range("K61").Value = 14709155
$ws.Range("M61").Value = -14708943

$ws.Range("H129").Value = 112400
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 112400
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 112400
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -122400

$ws.Range("H130").Value = 113000
$ws.Range("J130").Value = 113000
$ws.Range("L130").Value = 113000
$ws.Range("N130").Value = -123040

$ws.Range("H132").Value = 7573.4707
$ws.Range("I132").Value = 3237.1538
$ws.Range("K132").Value = 9711.4614
$ws.Range("M132").Value = -7181.4614

$ws.Range("H136").Value = 13188361
$ws.Range("I136").Value = 14709155
$ws.Range("K136").Value = 44127465
$ws.Range("M136").Value = -44124915

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 21922.75
$ws.Range("I26").Value = 10768.857
$ws.Range("K26").Value = 10768.857
$ws.Range("M26").Value = -10476.857

$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()

$ws.Range("H129").Value = 118434.25
$ws.Range("J129").Value = 118434.25
$ws.Range("L129").Value = 118434.25
$ws.Range("N129").Value = -128434.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 131999.5
$ws.Range("J20").Value = 131999.5
$ws.Range("L20").Value = 131999.5
$ws.Range("N20").Value = -132471.5

$ws.Range("H30").Value = 131999.5
$ws.Range("J30").Value = 131999.5
$ws.Range("L30").Value = 131999.5
$ws.Range("N30").Value = -132181.5

$ws.Range("H51").Value = 38000
$ws.Range("J51").Value = 38000
$ws.Range("L51").Value = 38000
$ws.Range("N51").Value = -39472

$ws.Range("H58").Value = 1481.9423
$ws.Range("I58").Value = 1175.2766
$ws.Range("K58").Value = 1175.2766
$ws.Range("M58").Value = -972.2765999999999

$ws.Range("H61").Value = 38000
$ws.Range("J61").Value = 38000
$ws.Range("L61").Value = 38000
$ws.Range("N61").Value = -38696

$ws.Range("H124").Value = 110000
$ws.Range("J124").Value = 110000
$ws.Range("L124").Value = 110000
$ws.Range("N124").Value = -114910

$ws.Range("H128").Value = 131999.5
$ws.Range("J128").Value = 131999.5
$ws.Range("L128").Value = 131999.5
$ws.Range("N128").Value = -141959.5

$ws.Range("H130").Value = 89000
$ws.Range("J130").Value = 89000
$ws.Range("L130").Value = 89000
$ws.Range("N130").Value = -99040

$ws.Range("H132").Value = 4431.4116
$ws.Range("I132").Value = 4357.5
$ws.Range("K132").Value = 13072.5
$ws.Range("M132").Value = -10542.5

$ws.Range("H134").Value = 360260.75
$ws.Range("I134").Value = 436751.44
$ws.Range("J134").Value = 8403.6
$ws.Range("K134").Value = 1310254.32
$ws.Range("L134").Value = 25210.8
$ws.Range("M134").Value = -1307719.32
$ws.Range("N134").Value = -30280.8

$ws.Range("H136").Value = 1481.9423
$ws.Range("I136").Value = 1175.2766
$ws.Range("K136").Value = 3525.8298
$ws.Range("M136").Value = -975.8297999999995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 8625621
$ws.Range("I4").Value = 12687682
$ws.Range("K4").Value = 38063046
$ws.Range("M4").Value = -38062934

$ws.Range("H5").Value = 2169
$ws.Range("J5").Value = 2200
$ws.Range("L5").Value = 6600
$ws.Range("N5").Value = -6824

$ws.Range("H50").Value = 421
$ws.Range("I50").Value = 102.5
$ws.Range("J50").Value = 448.69565
$ws.Range("K50").Value = 307.5
$ws.Range("L50").Value = 1346.08695
$ws.Range("M50").Value = 173.5
$ws.Range("N50").Value = -2308.08695

$ws.Range("H53").Value = 421
$ws.Range("I53").Value = 102.5
$ws.Range("J53").Value = 448.69565
$ws.Range("K53").Value = 307.5
$ws.Range("L53").Value = 1346.08695
$ws.Range("M53").Value = 173.5
$ws.Range("N53").Value = -2308.08695

$ws.Range("H58").Value = 55558056
$ws.Range("I58").Value = 1665.3334
$ws.Range("J58").Value = 111114450
$ws.Range("K58").Value = 4996.0002
$ws.Range("L58").Value = 333343350
$ws.Range("M58").Value = -4868.0002
$ws.Range("N58").Value = -333343606

$ws.Range("H129").Value = 17597026
$ws.Range("I129").Value = 3926.8
$ws.Range("J129").Value = 37144916
$ws.Range("K129").Value = 11780.4
$ws.Range("L129").Value = 111434748
$ws.Range("M129").Value = -6780.400000000001
$ws.Range("N129").Value = -111444748

$ws.Range("H130").Value = 4252.9414
$ws.Range("I130").Value = 1650
$ws.Range("K130").Value = 4950
$ws.Range("M130").Value = 70

$ws.Range("H131").Value = 6552.359
$ws.Range("I131").Value = 6623.9165
$ws.Range("J131").Value = 6520.5557
$ws.Range("K131").Value = 19871.7495
$ws.Range("L131").Value = 19561.6671
$ws.Range("M131").Value = -14831.7495
$ws.Range("N131").Value = -29641.6671

$ws.Range("H132").Value = 1650.3
$ws.Range("J132").Value = 2000.5
$ws.Range("L132").Value = 18004.5
$ws.Range("N132").Value = -23064.5

$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("M133").ClearContents()
$ws.Range("N133").ClearContents()

$ws.Range("H134").Value = 6354.933
$ws.Range("J134").Value = 8312.619000000001
$ws.Range("L134").Value = 24937.857
$ws.Range("N134").Value = -35077.857

$ws.Range("H135").Value = 2169
$ws.Range("J135").Value = 2200
$ws.Range("L135").Value = 19800
$ws.Range("N135").Value = -24870

$ws.Range("H136").Value = 5305.4116
$ws.Range("I136").Value = 5324.6875
$ws.Range("J136").Value = 4997
$ws.Range("K136").Value = 15974.0625
$ws.Range("L136").Value = 14991
$ws.Range("M136").Value = -10874.0625
$ws.Range("N136").Value = -25191

$ws.Range("H137").Value = 6700.5
$ws.Range("I137").Value = 4603.5
$ws.Range("J137").Value = 8098.5
$ws.Range("K137").Value = 13810.5
$ws.Range("L137").Value = 24295.5
$ws.Range("M137").Value = -8710.5
$ws.Range("N137").Value = -34495.5

$ws.Range("H138").Value = 2574.4
$ws.Range("J138").Value = 4000
$ws.Range("L138").Value = 12000
$ws.Range("N138").Value = -22280

$ws.Range("H139").Value = 2259.7646
$ws.Range("I139").Value = 1849.2858
$ws.Range("J139").Value = 2922.8462
$ws.Range("K139").Value = 5547.857400000001
$ws.Range("L139").Value = 8768.5386
$ws.Range("M139").Value = -407.8574000000008
$ws.Range("N139").Value = -19048.5386

$ws.Range("H140").Value = 752536
$ws.Range("I140").Value = 1002381.3
$ws.Range("K140").Value = 3007143.9
$ws.Range("M140").Value = -3001963.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()

$ws.Range("H14").Value = 1673934
$ws.Range("I14").Value = 640.8
$ws.Range("J14").Value = 2869143.5
$ws.Range("K14").Value = 640.8
$ws.Range("L14").Value = 2869143.5
$ws.Range("M14").Value = -472.8
$ws.Range("N14").Value = -2869479.5

$ws.Range("H47").Value = 30000
$ws.Range("J47").Value = 30000
$ws.Range("L47").Value = 30000
$ws.Range("N47").Value = -31136

$ws.Range("H102").Value = 1429.1897
$ws.Range("I102").Value = 1403.2858
$ws.Range("J102").Value = 1570.2222
$ws.Range("K102").Value = 1403.2858
$ws.Range("L102").Value = 1570.2222
$ws.Range("M102").Value = 218.7141999999999
$ws.Range("N102").Value = -4814.2222

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 30589.059
$ws.Range("I20").Value = 19999.715
$ws.Range("J20").Value = 80006
$ws.Range("K20").Value = 19999.715
$ws.Range("L20").Value = 80006
$ws.Range("M20").Value = -19773.715
$ws.Range("N20").Value = -80458

$ws.Range("H43").Value = 2521833
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 2521833
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 2521833
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -2522219

$ws.Range("H130").Value = 28057
$ws.Range("J130").Value = 28057
$ws.Range("L130").Value = 28057
$ws.Range("N130").Value = -38097

$ws.Range("H136").Value = 158279.75
$ws.Range("I136").Value = 17705.5
$ws.Range("K136").Value = 53116.5
$ws.Range("M136").Value = -50566.5

$ws.Range("H139").Value = 20000000
$ws.Range("I139").Value = 20000000
$ws.Range("K139").Value = 20000000
$ws.Range("M139").Value = -19994860

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 7000
$ws.Range("J14").Value = 7000
$ws.Range("L14").Value = 7000
$ws.Range("N14").Value = -7336

$ws.Range("H132").Value = 2585.26
$ws.Range("I132").Value = 2337.279
$ws.Range("J132").Value = 4108.5713
$ws.Range("K132").Value = 7011.837
$ws.Range("L132").Value = 12325.7139
$ws.Range("M132").Value = -4481.837
$ws.Range("N132").Value = -17385.7139

$ws.Range("H136").Value = 1935
$ws.Range("I136").Value = 1365.7727
$ws.Range("J136").Value = 3500.375
$ws.Range("K136").Value = 4097.3181
$ws.Range("L136").Value = 10501.125
$ws.Range("M136").Value = -1547.3181
$ws.Range("N136").Value = -15601.125
